$d = $word.ActiveDocument
$x = $d.Content.WordOpenXML

# ---- locate the document.xml part payload ----
$docStart = $x.IndexOf("<w:document ")
$docEndTag = "</w:document>"
$docEnd = $x.IndexOf($docEndTag) + $docEndTag.Length
$before = $x.Substring(0, $docStart)
$docXml = $x.Substring($docStart, $docEnd - $docStart)
$after = $x.Substring($docEnd)

# ---- paragraph 1: drop hyperlink + stray space run, add numPr + new runs ----
$p1Old = '<w:p w14:paraId="1240A510" w14:textId="387A6681" w:rsidR="004703A2" w:rsidRDefault="006B25B4" w:rsidP="00F55B81"><w:pPr><w:pStyle w:val="ListParagraph"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10" w:cs="CMR10"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:hyperlink r:id="rId5" w:history="1"><w:r w:rsidRPr="00B46D5C"><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10" w:cs="CMR10"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr><w:t>https://dakofa.dk/sektordatabase/</w:t></w:r></w:hyperlink><w:r><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10" w:cs="CMR10"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$p1New = '<w:p w14:paraId="1240A510" w14:textId="387A6681" w:rsidR="004703A2" w:rsidRDefault="006B25B4" w:rsidP="00F55B81"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10" w:cs="CMR10"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>You need to review the t</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">wo strands of the literature which are in the folder named social, in particular the papers highlighted in green in the word file called summaries. </w:t></w:r></w:p>'
if (-not $docXml.Contains($p1Old)) { throw "p1 anchor not found" }
$docXml = $docXml.Replace($p1Old, $p1New)

# ---- paragraph 2: add numPr + new run, then append six brand-new list paragraphs ----
$p2Old = '<w:p w14:paraId="55781A97" w14:textId="77777777" w:rsidR="004703A2" w:rsidRPr="00D55754" w:rsidRDefault="004703A2" w:rsidP="00F55B81"><w:pPr><w:pStyle w:val="ListParagraph"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10" w:cs="CMR10"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
$p2New = '<w:p w14:paraId="55781A97" w14:textId="77777777" w:rsidR="004703A2" w:rsidRPr="00D55754" w:rsidRDefault="004703A2" w:rsidP="00F55B81"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10" w:cs="CMR10"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">You need to review the immigration related situation in Denmark. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10" w:cs="CMR10"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">You need to review the literature on the energy sector in Denmark. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10" w:cs="CMR10"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">You need to review more in detail how the system of subsidies to the renewable energies works. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10" w:cs="CMR10"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">You need to review the chemical properties of dioxins </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10" w:cs="CMR10"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">You need to review the Danish situation in terms of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>pfoa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="CMR10" w:hAnsi="CMR10" w:cs="CMR10"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">See whether the iris (the risk assessment model used by the us </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>epa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">) has better estimates for the number of cancer cases associated with different toxic substances. </w:t></w:r></w:p>'
if (-not $docXml.Contains($p2Old)) { throw "p2 anchor not found" }
$docXml = $docXml.Replace($p2Old, $p2New)

$x = $before + $docXml + $after

# ---- locate the numbering.xml part payload ----
$numStart = $x.IndexOf("<w:numbering ")
$numEndTag = "</w:numbering>"
$numEnd = $x.IndexOf($numEndTag) + $numEndTag.Length
$nbefore = $x.Substring(0, $numStart)
$numXml = $x.Substring($numStart, $numEnd - $numStart)
$nafter = $x.Substring($numEnd)

# ---- renumber the old abstractNumId=3..9 definitions to 4..10 FIRST (descending order, no
#      collisions, and done before the new id-3 block exists so it cannot be touched twice) ----
$marker3 = '<w:abstractNum w:abstractNumId="3" w15:restartNumberingAfterBreak="0">'
if (-not $numXml.Contains($marker3)) { throw "abstractNum 3 marker not found" }
$numXml = $numXml.Replace('<w:abstractNum w:abstractNumId="9" w15:restartNumberingAfterBreak="0">', '<w:abstractNum w:abstractNumId="10" w15:restartNumberingAfterBreak="0">')
$numXml = $numXml.Replace('<w:abstractNum w:abstractNumId="8" w15:restartNumberingAfterBreak="0">', '<w:abstractNum w:abstractNumId="9" w15:restartNumberingAfterBreak="0">')
$numXml = $numXml.Replace('<w:abstractNum w:abstractNumId="7" w15:restartNumberingAfterBreak="0">', '<w:abstractNum w:abstractNumId="8" w15:restartNumberingAfterBreak="0">')
$numXml = $numXml.Replace('<w:abstractNum w:abstractNumId="6" w15:restartNumberingAfterBreak="0">', '<w:abstractNum w:abstractNumId="7" w15:restartNumberingAfterBreak="0">')
$numXml = $numXml.Replace('<w:abstractNum w:abstractNumId="5" w15:restartNumberingAfterBreak="0">', '<w:abstractNum w:abstractNumId="6" w15:restartNumberingAfterBreak="0">')
$numXml = $numXml.Replace('<w:abstractNum w:abstractNumId="4" w15:restartNumberingAfterBreak="0">', '<w:abstractNum w:abstractNumId="5" w15:restartNumberingAfterBreak="0">')
$numXml = $numXml.Replace($marker3, '<w:abstractNum w:abstractNumId="4" w15:restartNumberingAfterBreak="0">')

# ---- now insert the brand-new abstractNum (id 3, nsid 24D3727B) right where the old id-3 block used to be ----
$newAbstractNum = '<w:abstractNum w:abstractNumId="3" w15:restartNumberingAfterBreak="0"><w:nsid w:val="24D3727B"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="758A9566"/><w:lvl w:ilvl="0" w:tplc="262E1888"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%1."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1080" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorBidi" w:hint="default"/><w:color w:val="auto"/><w:sz w:val="22"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="04090019" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%2."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1800" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="2" w:tplc="0409001B" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%3."/><w:lvlJc w:val="right"/><w:pPr><w:ind w:left="2520" w:hanging="180"/></w:pPr></w:lvl><w:lvl w:ilvl="3" w:tplc="0409000F" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%4."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3240" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="4" w:tplc="04090019" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%5."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3960" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="5" w:tplc="0409001B" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%6."/><w:lvlJc w:val="right"/><w:pPr><w:ind w:left="4680" w:hanging="180"/></w:pPr></w:lvl><w:lvl w:ilvl="6" w:tplc="0409000F" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%7."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5400" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="7" w:tplc="04090019" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%8."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6120" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="8" w:tplc="0409001B" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%9."/><w:lvlJc w:val="right"/><w:pPr><w:ind w:left="6840" w:hanging="180"/></w:pPr></w:lvl></w:abstractNum>'
$oldBecame4 = '<w:abstractNum w:abstractNumId="4" w15:restartNumberingAfterBreak="0">'
$numXml = $numXml.Replace($oldBecame4, $newAbstractNum + $oldBecame4)

# ---- repoint every <w:num> entry at its (possibly shifted) abstractNumId ----
$numXml = $numXml.Replace('<w:num w:numId="1"><w:abstractNumId w:val="9"/></w:num>', '<w:num w:numId="1"><w:abstractNumId w:val="10"/></w:num>')
$numXml = $numXml.Replace('<w:num w:numId="2"><w:abstractNumId w:val="8"/></w:num>', '<w:num w:numId="2"><w:abstractNumId w:val="9"/></w:num>')
$numXml = $numXml.Replace('<w:num w:numId="3"><w:abstractNumId w:val="0"/></w:num>', '<w:num w:numId="3"><w:abstractNumId w:val="0"/></w:num>')
$numXml = $numXml.Replace('<w:num w:numId="4"><w:abstractNumId w:val="7"/></w:num>', '<w:num w:numId="4"><w:abstractNumId w:val="8"/></w:num>')
$numXml = $numXml.Replace('<w:num w:numId="5"><w:abstractNumId w:val="6"/></w:num>', '<w:num w:numId="5"><w:abstractNumId w:val="7"/></w:num>')
$numXml = $numXml.Replace('<w:num w:numId="6"><w:abstractNumId w:val="4"/></w:num>', '<w:num w:numId="6"><w:abstractNumId w:val="5"/></w:num>')
$numXml = $numXml.Replace('<w:num w:numId="7"><w:abstractNumId w:val="3"/></w:num>', '<w:num w:numId="7"><w:abstractNumId w:val="4"/></w:num>')
$numXml = $numXml.Replace('<w:num w:numId="8"><w:abstractNumId w:val="1"/></w:num>', '<w:num w:numId="8"><w:abstractNumId w:val="1"/></w:num>')
$numXml = $numXml.Replace('<w:num w:numId="9"><w:abstractNumId w:val="2"/></w:num>', '<w:num w:numId="9"><w:abstractNumId w:val="2"/></w:num>')
$numXml = $numXml.Replace('<w:num w:numId="10"><w:abstractNumId w:val="5"/></w:num>', '<w:num w:numId="10"><w:abstractNumId w:val="6"/></w:num>')

# ---- register the new list (numId 11 -> abstractNumId 3) used by the paragraphs above ----
$numXml = $numXml.Replace("</w:numbering>", '<w:num w:numId="11"><w:abstractNumId w:val="3"/></w:num></w:numbering>')

$x = $nbefore + $numXml + $nafter

$d.Content.InsertXML($x)
Write-Output "done"
